$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04332866666666666
$ws.Range("H2").Value = 0.129986
$ws.Range("I2").Value = 0.2246397599897691
$ws.Range("J2").Value = 0.2246397599897691
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.404869
$ws.Range("N2").Value = 4.214607
$ws.Range("O2").Value = 0.6692718564235921
$ws.Range("P2").Value = 0.6692718564235923
$ws.Range("Q2").Value = 0.06087110061133332
$ws.Range("R2").Value = 0.547839905502
$ws.Range("S2").Value = 0.150345069194903
$ws.Range("T2").Value = 0.150345069194903
$ws.Range("G3").Value = 0.04332866666666666
$ws.Range("H3").Value = 0.129986
$ws.Range("I3").Value = 0.2246397599897691
$ws.Range("J3").Value = 0.2246397599897691
$ws.Range("O3").Value = 0.3150411080808892
$ws.Range("P3").Value = 0.3150411080808893
$ws.Range("Q3").Value = 0.02865337725266666
$ws.Range("R3").Value = 0.257880395274
$ws.Range("S3").Value = 0.07077075890620187
$ws.Range("T3").Value = 0.07077075890620188
$ws.Range("G4").Value = 0.04332866666666666
$ws.Range("H4").Value = 0.129986
$ws.Range("I4").Value = 0.2246397599897691
$ws.Range("J4").Value = 0.2246397599897691
$ws.Range("M4").Value = 0.03292866666666667
$ws.Range("N4").Value = 0.098786
$ws.Range("O4").Value = 0.01568703549551856
$ws.Range("P4").Value = 0.01568703549551856
$ws.Range("Q4").Value = 0.001426755221777778
$ws.Range("R4").Value = 0.012840796996
$ws.Range("S4").Value = 0.003523931888664278
$ws.Range("T4").Value = 0.003523931888664279
$ws.Range("I5").Value = 0.5955530362469368
$ws.Range("J5").Value = 0.5955530362469369
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.404869
$ws.Range("N5").Value = 4.214607
$ws.Range("O5").Value = 0.6692718564235921
$ws.Range("P5").Value = 0.6692718564235923
$ws.Range("Q5").Value = 0.1613782386093333
$ws.Range("R5").Value = 1.452404147484
$ws.Range("S5").Value = 0.3985868861676942
$ws.Range("T5").Value = 0.3985868861676944
$ws.Range("I6").Value = 0.5955530362469368
$ws.Range("J6").Value = 0.5955530362469369
$ws.Range("O6").Value = 0.3150411080808892
$ws.Range("P6").Value = 0.3150411080808893
$ws.Range("S6").Value = 0.1876236884601729
$ws.Range("T6").Value = 0.187623688460173
$ws.Range("I7").Value = 0.5955530362469368
$ws.Range("J7").Value = 0.5955530362469369
$ws.Range("M7").Value = 0.03292866666666667
$ws.Range("N7").Value = 0.098786
$ws.Range("O7").Value = 0.01568703549551856
$ws.Range("P7").Value = 0.01568703549551856
$ws.Range("Q7").Value = 0.003782537892444445
$ws.Range("R7").Value = 0.034042841032
$ws.Range("S7").Value = 0.009342461619069549
$ws.Range("T7").Value = 0.00934246161906955
$ws.Range("G8").Value = 0.03468133333333333
$ws.Range("H8").Value = 0.104044
$ws.Range("I8").Value = 0.1798072037632941
$ws.Range("J8").Value = 0.1798072037632941
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.404869
$ws.Range("N8").Value = 4.214607
$ws.Range("O8").Value = 0.6692718564235921
$ws.Range("P8").Value = 0.6692718564235923
$ws.Range("Q8").Value = 0.04872273007866667
$ws.Range("R8").Value = 0.438504570708
$ws.Range("S8").Value = 0.1203399010609949
$ws.Range("T8").Value = 0.120339901060995
$ws.Range("G9").Value = 0.03468133333333333
$ws.Range("H9").Value = 0.104044
$ws.Range("I9").Value = 0.1798072037632941
$ws.Range("J9").Value = 0.1798072037632941
$ws.Range("O9").Value = 0.3150411080808892
$ws.Range("P9").Value = 0.3150411080808893
$ws.Range("Q9").Value = 0.02293486977733333
$ws.Range("R9").Value = 0.206413827996
$ws.Range("S9").Value = 0.05664666071451439
$ws.Range("T9").Value = 0.0566466607145144
$ws.Range("G10").Value = 0.03468133333333333
$ws.Range("H10").Value = 0.104044
$ws.Range("I10").Value = 0.1798072037632941
$ws.Range("J10").Value = 0.1798072037632941
$ws.Range("M10").Value = 0.03292866666666667
$ws.Range("N10").Value = 0.098786
$ws.Range("O10").Value = 0.01568703549551856
$ws.Range("P10").Value = 0.01568703549551856
$ws.Range("Q10").Value = 0.001142010064888889
$ws.Range("R10").Value = 0.010278090584
$ws.Range("S10").Value = 0.002820641987784732
$ws.Range("T10").Value = 0.002820641987784732
